# Applies the diff: rename column header Gen->MaxFES, change MaxFES values,
# fill in previously-blank X-column values, drop the "Run 50" run column
# (shifting Mean into its place with recomputed values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header A1: "Gen" -> "MaxFES" ---
$ws.Range("A1").Value = "MaxFES"

# --- 2. Column A (MaxFES) values for rows 3-14 ---
$newA = @{
    3  = 0.001
    4  = 0.01
    5  = 0.1
    6  = 0.2
    7  = 0.3
    8  = 0.4
    9  = 0.5
    10 = 0.6
    11 = 0.7
    12 = 0.8
    13 = 0.9
    14 = 1
}
foreach ($row in $newA.Keys) {
    $ws.Cells.Item($row, 1).Value = $newA[$row]
}

# --- 3. Column X (Run 22) values for rows 4-14, previously blank ---
$newX = @{
    4  = 9.2942009
    5  = 5.48615013
    6  = 3.54259594
    7  = 0.82218494
    8  = 0.25721916
    9  = 0.0872908
    10 = 0.01391343
    11 = 0.0024497
    12 = 0.00030457
    13 = 0.00003906
    14 = 0.00000001
}
foreach ($row in $newX.Keys) {
    $ws.Cells.Item($row, 24).Value = $newX[$row]
}

# --- 4. Drop the "Run 50" column (AZ). This removes that run's data and
#        shifts the "Mean" column (formerly BA) into AZ. ---
$ws.Columns("AZ").Delete()

# --- 5. Update the (now-shifted) Mean column AZ with recomputed values ---
$newMean = @{
    2  = 13.75202703
    3  = 12.72643367
    4  = 9.50078501
    5  = 4.91997835
    6  = 3.01878766
    7  = 1.88201447
    8  = 1.30131162
    9  = 0.94958964
    10 = 0.74292487
    11 = 0.6514306
    12 = 0.59757852
    13 = 0.56640335
    14 = 0.55149848
}
foreach ($row in $newMean.Keys) {
    $ws.Cells.Item($row, 52).Value = $newMean[$row]
}
